$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 17 and 18 by copying the formatting of the last data row (16),
# then overwrite with the new record's values - this keeps number formats
# (date column G, centered L/M columns, etc.) consistent with the rest of
# the table.
$ws.Range("B16:M16").Copy()
$ws.Range("B17:M17").PasteSpecial(-4122)
$ws.Range("B16:M16").Copy()
$ws.Range("B18:M18").PasteSpecial(-4122)

# --- Row 17: new item 10011 ---
$ws.Range("B17").Value = 10011
$ws.Range("C17").Value = "测试"
$ws.Range("D17").Value = "测试223"
$ws.Range("E17").Value = 1222
$ws.Range("F17").Value = 10011
$ws.Range("G17").Value = 45930.625
$ws.Range("H17").Value = $true
$ws.Range("I17").Value = "RED"
$ws.Range("J17").Value = "10011,9"
$ws.Range("K17").Value = "10001,2;10002,12"
$ws.Range("L17").Value = 10001
$ws.Range("M17").Value = 12

# --- Row 18: new item 10012 (desc typed with a leading apostrophe to force text) ---
$ws.Range("B18").Value = 10012
$ws.Range("C18").Value = "测试123"
$ws.Range("D18").Value = "'121121213312"
$ws.Range("E18").Value = 1222
$ws.Range("F18").Value = 10012
$ws.Range("G18").Value = 45930.625
$ws.Range("H18").Value = $true
$ws.Range("I18").Value = "RED"
$ws.Range("J18").Value = "10011,9"
$ws.Range("K18").Value = "10001,2;10002,12"
$ws.Range("L18").Value = 10001
$ws.Range("M18").Value = 12

# --- View state: restore the selection left by the editing session ---
$ws.Range("H27").Select() | Out-Null
